$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "53.791.05"
Set-TextValue $ws.Range("E2") "  -9.28%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.409.59"
Set-TextValue $ws.Range("E3") "  -12.31%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.03%  "

# Row 5
Set-TextValue $ws.Range("D5") "462.13"
Set-TextValue $ws.Range("E5") "  -9.53%  "

# Row 6
Set-TextValue $ws.Range("D6") "130.68"
Set-TextValue $ws.Range("E6") "  -8.88%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.996"
Set-TextValue $ws.Range("E7") "  -0.04%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.485"
Set-TextValue $ws.Range("E8") "  -10.03%  "

# Row 9
Set-TextValue $ws.Range("D9") "2.429.49"
Set-TextValue $ws.Range("E9") "  -12.07%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0950"
Set-TextValue $ws.Range("E10") "  -10.20%  "

# Row 11
Set-TextValue $ws.Range("D11") "5.32"
Set-TextValue $ws.Range("E11") "  -13.09%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.313"
Set-TextValue $ws.Range("E12") "  -11.39%  "

# Row 13
Set-TextValue $ws.Range("E13") "  -4.51%  "

# Row 14
Set-TextValue $ws.Range("D14") "2.825.60"
Set-TextValue $ws.Range("E14") "  -12.54%  "

# Row 15
Set-TextValue $ws.Range("D15") "53.790.26"
Set-TextValue $ws.Range("E15") "  -9.23%  "

# Row 16
Set-TextValue $ws.Range("D16") "19.65"
Set-TextValue $ws.Range("E16") "  -10.90%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.0000131"
Set-TextValue $ws.Range("E17") "  -4.49%  "

# Row 18
Set-TextValue $ws.Range("D18") "2.411.25"
Set-TextValue $ws.Range("E18") "  -13.06%  "

# Row 19
Set-TextValue $ws.Range("D19") "4.15"
Set-TextValue $ws.Range("E19") "  -13.73%  "

# Row 20
Set-TextValue $ws.Range("D20") "306.63"
Set-TextValue $ws.Range("E20") "  -12.18%  "

# Row 21
Set-TextValue $ws.Range("D21") "9.37"
Set-TextValue $ws.Range("E21") "  -16.16%  "

# Row 22
Set-TextValue $ws.Range("D22") "1.00"
Set-TextValue $ws.Range("E22") "  +0.18%  "

# Row 23
Set-TextValue $ws.Range("E23") "  +1.00%  "

# Row 24
Set-TextValue $ws.Range("D24") "5.35"
Set-TextValue $ws.Range("E24") "  -15.38%  "

# Row 25
Set-TextValue $ws.Range("D25") "55.78"
Set-TextValue $ws.Range("E25") "  -12.16%  "

# Row 26
Set-TextValue $ws.Range("E26") "  +0.86%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.384"
Set-TextValue $ws.Range("E27") "  -10.88%  "

# Row 28
Set-TextValue $ws.Range("B28") "Kaspa"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D28") "0.153"
Set-TextValue $ws.Range("E28") "  -12.50%  "

# Row 29
Set-TextValue $ws.Range("B29") "WrappedeETH"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D29") "2.476.17"
Set-TextValue $ws.Range("E29") "  -13.35%  "

# Row 30
Set-TextValue $ws.Range("D30") "7.07"
Set-TextValue $ws.Range("E30") "  -7.07%  "

# Row 31
Set-TextValue $ws.Range("E31") "  -0.07%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.0₃0718"
Set-TextValue $ws.Range("E32") "  -15.65%  "

# Row 33
Set-TextValue $ws.Range("D33") "145.52"
Set-TextValue $ws.Range("E33") "  -3.05%  "

# Row 34
Set-TextValue $ws.Range("D34") "17.68"
Set-TextValue $ws.Range("E34") "  -8.97%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.42"
Set-TextValue $ws.Range("E35") "  -12.47%  "

# Row 36
Set-TextValue $ws.Range("D36") "4.98"
Set-TextValue $ws.Range("E36") "  -8.99%  "

# Row 37
Set-TextValue $ws.Range("D37") "3.51"
Set-TextValue $ws.Range("E37") "  -17.81%  "

# Row 38
Set-TextValue $ws.Range("D38") "1.06"
Set-TextValue $ws.Range("E38") "  -8.02%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.798"
Set-TextValue $ws.Range("E39") "  -17.60%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.995"
Set-TextValue $ws.Range("E40") "  +0.03%  "

# Row 41
Set-TextValue $ws.Range("D41") "32.93"
Set-TextValue $ws.Range("E41") "  -9.25%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.594"
Set-TextValue $ws.Range("E42") "  -2.74%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.0524"
Set-TextValue $ws.Range("E43") "  -7.12%  "

# Row 44
Set-TextValue $ws.Range("D44") "3.25"
Set-TextValue $ws.Range("E44") "  -8.95%  "

# Row 45
Set-TextValue $ws.Range("D45") "10.10"
Set-TextValue $ws.Range("E45") "  -2.66%  "

# Row 46
Set-TextValue $ws.Range("D46") "1.23"
Set-TextValue $ws.Range("E46") "  -13.06%  "

# Row 47
Set-TextValue $ws.Range("D47") "1.930.64"
Set-TextValue $ws.Range("E47") "  -12.34%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.0217"
Set-TextValue $ws.Range("E48") "  -5.17%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0866"
Set-TextValue $ws.Range("E49") "  -2.98%  "

# Row 50
Set-TextValue $ws.Range("B50") "RenderToken"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D50") "4.24"
Set-TextValue $ws.Range("E50") "  -11.88%  "

# Row 51
Set-TextValue $ws.Range("B51") "EnergySwap"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "16.54"
Set-TextValue $ws.Range("E51") "  -13.35%  "
